# Avances Etiquetado Roboflow - weekly update for 6/10/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in this week's row (row 33) of the tracking table
$ws.Range("D33").Value = 45936   # Fecha (serial date -> 06/10/2025)
$ws.Range("E33").Value = 75      # Imagenes sin etiquetar
$ws.Range("F33").Value = 443     # Imagenes etiquetadas sin revisar
$ws.Range("G33").Value = 0       # Imagenes rechazadas
$ws.Range("H33").Value = 0       # Imagenes etiquetadas y revisadas, faltando de subir
$ws.Range("I33").Value = 1012    # Imagenes etiquetadas, revisadas y subidas
$ws.Range("J33").Value = "N/A"   # Notas

# Leave the view/selection where the author left it when saving
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J34").Select()
